$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15: repeat header row (Nombre / Tipo / Descripción) ---
$ws.Range("A15").Value = "Nombre"
$ws.Range("B15").Value = "Tipo"
$ws.Range("C15").Value = "Descripción"

# Copy header formatting (fill + border + centered) from row 5 to row 15
$ws.Range("A5:C5").Copy()
$ws.Range("A15:C15").PasteSpecial(-4122)

# --- Row 16: new actor "Comprador" ---
$ws.Range("A16").Value = "Comprador"
$ws.Range("B16").Value = "Cliente"
$ws.Range("C16").Value = "Persona que desea adquirir una montura nueva o unos lentes o reparar sus lentes o montura o requiere de un examen visual"
$ws.Range("D16").Value = "> Esto es en caso de que sea solo con un cliente (que no creo)"

# Copy the data-row formatting (border + wrap + vertical center) from row 6 to A16:C16
$ws.Range("A6:C6").Copy()
$ws.Range("A16:C16").PasteSpecial(-4122)

# D16 gets the "blank separator" formatting (wrap + vcenter, no outer border) as a base,
# then a left-only thin border is added to visually separate it from column C
$ws.Range("A11").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D16").Borders(7).LineStyle = 1

$ws.Rows("16").RowHeight = 90

# --- View state ---
$ws.Range("E11").Select()
$excel.ActiveWindow.ScrollRow = 7
